$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace Matheus Diniz's entry with a Doris entry ---
$ws.Range("A2").Value = "Dóris Andressa Moura Luvizute"

$ws.Range("B2").Value = "dorisluvizute@gmail.com"
$ws.Range("B2").Style = "Hiperlink"
$ws.Range("B2").Font.Underline = 2

# C2 ("ENVIADO") is unchanged.

# New empty styled cell E2 (underlined black font, Normal style / xfId 0)
$ws.Range("E2").Font.Underline = 2

# --- Row 3: Doris's email row keeps its text but gains underline styling ---
$ws.Range("C3").Font.Underline = 2

# --- Row 4: clear out the old Doris/ENVIADO row, leaving only a styled blank B4 ---
$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()

# --- Move the mailto hyperlink from B2 to B3 (same target address) ---
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:matheusinhodinizinho@gmail.com")

# Re-assert B3's explicit formatting so it matches a plain underlined
# Hiperlink-style cell (Hyperlinks.Add leaves a redundant applyFont flag).
$ws.Range("B3").Font.Underline = 2
$ws.Range("B3").Style = "Hiperlink"

# --- Update the selection shown when the sheet is opened ---
$ws.Range("C3").Select() | Out-Null
